$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.281.53'
$ws.Range('E2').Value = '  -0.55%  '
$ws.Range('D3').Value = '2.066.30'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.83'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.620'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.34%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '56.94'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.383'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0762'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.12%  '
$ws.Range('E11').Value = '  +0.53%  '
$ws.Range('D12').Value = '2.370.15'
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.62'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.79'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.60%  '
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.14'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.94%  '
$ws.Range('D17').Value = '2.064.94'
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('D18').Value = '37.208.45'
$ws.Range('E18').Value = '  -0.96%  '
$ws.Range('E19').Value = '  +4.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.49'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.41%  '
$ws.Range('D21').Value = '0.0₃0812'
$ws.Range('E21').Value = '  +0.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '225.68'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.35%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('E24').Value = '  +0.97%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.40'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.44%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.42'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.78'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.43'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.95'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.58%  '
$ws.Range('E30').Value = '  -3.13%  '
$ws.Range('E31').Value = '  -1.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.49'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0617'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.93%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.60'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.50'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.89%  '
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('E37').Value = '  -1.75%  '
$ws.Range('E38').Value = '  -3.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.66'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.42%  '
$ws.Range('E40').Value = '  -1.47%  '
$ws.Range('D41').Value = '1.477.00'
$ws.Range('E41').Value = '  +0.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '96.07'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0932'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.11%  '
$ws.Range('E44').Value = '  +3.53%  '
$ws.Range('E45').Value = '  -0.16%  '
$ws.Range('E46').Value = '  -4.94%  '
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '15.22'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.18'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.97'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.32%  '
$ws.Range('D51').Value = '2.258.58'
$ws.Range('E51').Value = '  -0.11%  '
